# Updates bootstrap statistics (columns H:AE) on the "fractions" sheet for rows 4-13,
# reflecting a recomputation of the hospital age-stratified parameters.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fractions")

# Row 4
$ws.Range("H4").Value = 0.981
$ws.Range("I4").Value = 0.981
$ws.Range("L4").Value = 0.9710000000000001

# Row 5
$ws.Range("J5").Value = 0.862
$ws.Range("K5").Value = 0.878
$ws.Range("L5").Value = 0.847
$ws.Range("P5").Value = 0.008999999999999999
$ws.Range("R5").Value = 0.005475000000000002
$ws.Range("Z5").Value = 0.09163486005089058
$ws.Range("AA5").Value = 0.09163486005089058
$ws.Range("AB5").Value = 0.07349568267674042
$ws.Range("AC5").Value = 0.1094890510948905
$ws.Range("AD5").Value = 0.04540437383948834
$ws.Range("AE5").Value = 0.1462462462462462

# Row 6
$ws.Range("H6").Value = 0.9098196392785571
$ws.Range("I6").Value = 0.9098196392785571
$ws.Range("J6").Value = 0.9038076152304609
$ws.Range("K6").Value = 0.9159369369369369
$ws.Range("L6").Value = 0.8938380765535074
$ws.Range("M6").Value = 0.9269269269269269
$ws.Range("N6").Value = 0.01601601601601602
$ws.Range("O6").Value = 0.01601601601601602
$ws.Range("P6").Value = 0.01302605210420842
$ws.Range("Q6").Value = 0.01806323493092343
$ws.Range("S6").Value = 0.02501313813813814
$ws.Range("T6").Value = 0.007667031763417305
$ws.Range("U6").Value = 0.007667031763417305
$ws.Range("V6").Value = 0.005488474204171241
$ws.Range("W6").Value = 0.009849523416922116
$ws.Range("X6").Value = 0.002709165751258436
$ws.Range("Y6").Value = 0.01416200815689141
$ws.Range("Z6").Value = 0.09375
$ws.Range("AA6").Value = 0.09375
$ws.Range("AB6").Value = 0.07486702127659574
$ws.Range("AC6").Value = 0.1168831168831169
$ws.Range("AD6").Value = 0.04040404040404041
$ws.Range("AE6").Value = 0.1547435897435897

# Row 7
$ws.Range("H7").Value = 0.897
$ws.Range("I7").Value = 0.897
$ws.Range("J7").Value = 0.89
$ws.Range("K7").Value = 0.9029029029029029
$ws.Range("L7").Value = 0.877
$ws.Range("M7").Value = 0.913913913913914
$ws.Range("N7").Value = 0.02655220883534137
$ws.Range("O7").Value = 0.02655220883534137
$ws.Range("P7").Value = 0.02302302302302302
$ws.Range("Q7").Value = 0.03006765788347005
$ws.Range("R7").Value = 0.01702511629866339
$ws.Range("S7").Value = 0.03707414829659319
$ws.Range("T7").Value = 0.01128031943726921
$ws.Range("U7").Value = 0.01128031943726921
$ws.Range("V7").Value = 0.008946046267481534
$ws.Range("W7").Value = 0.01434878587196468
$ws.Range("X7").Value = 0.005521660385642677
$ws.Range("Y7").Value = 0.01882720061613327
$ws.Range("Z7").Value = 0.160188679245283
$ws.Range("AA7").Value = 0.160188679245283
$ws.Range("AB7").Value = 0.134020618556701
$ws.Range("AC7").Value = 0.1834862385321101
$ws.Range("AD7").Value = 0.08919289250353607
$ws.Range("AE7").Value = 0.2319555868000991

# Row 8
$ws.Range("H8").Value = 0.8582168674698796
$ws.Range("I8").Value = 0.8582168674698796
$ws.Range("J8").Value = 0.85
$ws.Range("K8").Value = 0.8656304585097979
$ws.Range("L8").Value = 0.8358358358358359
$ws.Range("M8").Value = 0.879765422393306
$ws.Range("N8").Value = 0.04104104104104104
$ws.Range("O8").Value = 0.04104104104104104
$ws.Range("P8").Value = 0.03702777777777778
$ws.Range("Q8").Value = 0.04604604604604605
$ws.Range("R8").Value = 0.02905811623246493
$ws.Range("S8").Value = 0.05413953109736242
$ws.Range("T8").Value = 0.02121397319258282
$ws.Range("U8").Value = 0.02121397319258282
$ws.Range("W8").Value = 0.02448819296645384
$ws.Range("X8").Value = 0.01282051282051282
$ws.Range("Y8").Value = 0.03119939320354708
$ws.Range("Z8").Value = 0.1620669965740388
$ws.Range("AA8").Value = 0.1620669965740388
$ws.Range("AB8").Value = 0.1397788258253375
$ws.Range("AC8").Value = 0.1824447622805287
$ws.Range("AD8").Value = 0.1017496054448609
$ws.Range("AE8").Value = 0.2222069269313363

# Row 9
$ws.Range("H9").Value = 0.8086172344689379
$ws.Range("I9").Value = 0.8086172344689379
$ws.Range("J9").Value = 0.8006042296072508
$ws.Range("L9").Value = 0.7845691382765531
$ws.Range("M9").Value = 0.8333333333333334
$ws.Range("N9").Value = 0.0802407221664995
$ws.Range("O9").Value = 0.0802407221664995
$ws.Range("Q9").Value = 0.08610765073690926
$ws.Range("R9").Value = 0.06516129048468049
$ws.Range("S9").Value = 0.09628886659979939
$ws.Range("T9").Value = 0.03625971502590673
$ws.Range("U9").Value = 0.03625971502590673
$ws.Range("V9").Value = 0.0322180916976456
$ws.Range("W9").Value = 0.04054482238892886
$ws.Range("X9").Value = 0.02454971410419314
$ws.Range("Y9").Value = 0.04915415951972556
$ws.Range("Z9").Value = 0.2661571000214179
$ws.Range("AA9").Value = 0.2661571000214179
$ws.Range("AB9").Value = 0.245048019207683
$ws.Range("AC9").Value = 0.2848993498374593
$ws.Range("AD9").Value = 0.2025279123657047
$ws.Range("AE9").Value = 0.3299763635439435

# Row 10
$ws.Range("H10").Value = 0.7598798798798798
$ws.Range("I10").Value = 0.7598798798798798
$ws.Range("J10").Value = 0.7492477432296891
$ws.Range("K10").Value = 0.7685370741482966
$ws.Range("L10").Value = 0.7298199554976663
$ws.Range("M10").Value = 0.7847557678170411
$ws.Range("P10").Value = 0.1571571571571572
$ws.Range("Q10").Value = 0.172043043043043
$ws.Range("R10").Value = 0.1412083175359728
$ws.Range("S10").Value = 0.1867020541082164
$ws.Range("T10").Value = 0.07572739345647712
$ws.Range("U10").Value = 0.07572739345647712
$ws.Range("V10").Value = 0.06948921398462683
$ws.Range("W10").Value = 0.08226437000313597
$ws.Range("X10").Value = 0.05760185723827375
$ws.Range("Y10").Value = 0.09440262699665479
$ws.Range("Z10").Value = 0.4444444444444444
$ws.Range("AA10").Value = 0.4444444444444444
$ws.Range("AB10").Value = 0.424308439329039
$ws.Range("AC10").Value = 0.4643647281921618
$ws.Range("AD10").Value = 0.3823839662447257
$ws.Range("AE10").Value = 0.5021948051948052

# Row 11
$ws.Range("H11").Value = 0.7838516908868667
$ws.Range("I11").Value = 0.7838516908868667
$ws.Range("J11").Value = 0.7755511022044088
$ws.Range("K11").Value = 0.7927927927927928
$ws.Range("L11").Value = 0.76024749498998
$ws.Range("M11").Value = 0.807431730931743
$ws.Range("N11").Value = 0.2665330661322645
$ws.Range("O11").Value = 0.2665330661322645
$ws.Range("P11").Value = 0.2575150300601202
$ws.Range("Q11").Value = 0.2750688188188188
$ws.Range("R11").Value = 0.2408420945996096
$ws.Range("S11").Value = 0.2934475808473804
$ws.Range("T11").Value = 0.1719543913614146
$ws.Range("U11").Value = 0.1719543913614146
$ws.Range("V11").Value = 0.1625761553520804
$ws.Range("W11").Value = 0.1825119493289804
$ws.Range("X11").Value = 0.1465863425831045
$ws.Range("Y11").Value = 0.1979448718463763
$ws.Range("Z11").Value = 0.6026785714285714
$ws.Range("AA11").Value = 0.6026785714285714
$ws.Range("AB11").Value = 0.5802678571428572
$ws.Range("AC11").Value = 0.6286067892503536
$ws.Range("AD11").Value = 0.5309557383470427

# Row 12
$ws.Range("K12").Value = 0.9237713139418254
$ws.Range("L12").Value = 0.8989272042854196
$ws.Range("M12").Value = 0.9327663462274374
$ws.Range("N12").Value = 0.4046184738955823
$ws.Range("O12").Value = 0.4046184738955823
$ws.Range("P12").Value = 0.3945323335383576
$ws.Range("Q12").Value = 0.41528815676258
$ws.Range("R12").Value = 0.3769393216080402
$ws.Range("S12").Value = 0.4348697394789579
$ws.Range("T12").Value = 0.3735667155862595
$ws.Range("U12").Value = 0.3735667155862595
$ws.Range("V12").Value = 0.3622988906693843
$ws.Range("W12").Value = 0.3837341285109658
$ws.Range("X12").Value = 0.3418275723630577
$ws.Range("Y12").Value = 0.4037115318674725
$ws.Range("Z12").Value = 0.7528089887640449
$ws.Range("AA12").Value = 0.7528089887640449
$ws.Range("AB12").Value = 0.72
$ws.Range("AC12").Value = 0.7857142857142857
$ws.Range("AD12").Value = 0.6535769230769231
$ws.Range("AE12").Value = 0.8395061728395061

# Row 13
$ws.Range("H13").Value = 0.8371746987951807
$ws.Range("I13").Value = 0.8371746987951807
$ws.Range("J13").Value = 0.8286573146292585
$ws.Range("K13").Value = 0.8454202970356853
$ws.Range("L13").Value = 0.8155341776833507
$ws.Range("M13").Value = 0.8622557612778274
$ws.Range("N13").Value = 0.2151452960498187
$ws.Range("O13").Value = 0.2151452960498187
$ws.Range("Q13").Value = 0.2238955823293173
$ws.Range("R13").Value = 0.1897735707121364
$ws.Range("S13").Value = 0.2396994861725155
$ws.Range("T13").Value = 0.1650543923265879
$ws.Range("U13").Value = 0.1650543923265879
$ws.Range("V13").Value = 0.1567004861551469
$ws.Range("W13").Value = 0.1732894258490518
$ws.Range("X13").Value = 0.1389973263165031
$ws.Range("Y13").Value = 0.1899597576084187
$ws.Range("Z13").Value = 0.4632509412342446
$ws.Range("AA13").Value = 0.4632509412342446
$ws.Range("AB13").Value = 0.440188679245283
$ws.Range("AC13").Value = 0.4924028822055138
$ws.Range("AD13").Value = 0.3944725028058361
$ws.Range("AE13").Value = 0.5548832271762207
